$wb = $excel.ActiveWorkbook

# 1. Rename "Sheet1" to "Data"
$wsData = $wb.Worksheets.Item("Sheet1")
$wsData.Name = "Data"

# 2. Rebuild the "Legend" sheet content as a 2-column table (Column1 / Column2)
$wsLegend = $wb.Worksheets.Item("Legend")

# Clear out the old single-column legend values first
$wsLegend.Range("A1:A6").Clear()

# New header row
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

# New data rows: numeric rating + effectiveness label (no leading digit anymore)
$wsLegend.Range("A2").Value = 0
$wsLegend.Range("B2").Value = " Not Effective "

$wsLegend.Range("A3").Value = 1
$wsLegend.Range("B3").Value = " Minimally Effective "

$wsLegend.Range("A4").Value = 2
$wsLegend.Range("B4").Value = " Somewhat Effective "

$wsLegend.Range("A5").Value = 3
$wsLegend.Range("B5").Value = " Moderately Effective "

$wsLegend.Range("A6").Value = 4
$wsLegend.Range("B6").Value = " Highly Effective "

$wsLegend.Range("A7").Value = 5
$wsLegend.Range("B7").Value = " Extremely Effective "

$wsLegend.Columns.Item(2).ColumnWidth = 11.28515625

# Turn the range into a native Excel Table ("Table1")
$tbl = $wsLegend.ListObjects.Add(1, $wsLegend.Range("A1:B7"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium2"

# 3. Make "Legend" the active/selected tab (was "Data" before)
$wsLegend.Activate()
$wsLegend.Range("A1:B7").Select()
